# 16.4.2.1 Number of seized/voluntarily surrendered firearms — update to add a
# 2020 column and drop the "Abducted"/"Lost" breakdown rows (now merged away).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ("Seized firearms") becomes the new last data row, so give it the
# plain bottom-border look that row 2 (the blank divider row) already uses,
# instead of the unbordered look it has today.
$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)

# --- Drop rows 6 and 7 ("Abducted" / "Lost"): that breakdown no longer ships.
$ws.Rows("6:7").Delete()

# --- Add the 2020 column (H), mirroring column G's formatting down each row.
$ws.Range("G2:G5").Copy()
$ws.Range("H2:H5").PasteSpecial(-4122)

# --- Updated 2019 figures and brand-new 2020 figures.
$ws.Range("G4").Value = 146
$ws.Range("G5").Value = 127
$ws.Range("H3").Value = 2020
$ws.Range("H4").Value = 158
$ws.Range("H5").Value = 397

$ws.Range("A1").Select()
